$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36-39: coin pairs reordered/updated (do these first, in ascending row order) ---
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'6.823"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02975"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'11.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.22%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2747"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.70%  "

# --- Price (D) and Volume(1h) (E) updates for all other rows ---
$ws.Range("D2").Value = "24.456.13"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.666.16"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'0.9937"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.05%  "
$ws.Range("D5").Value = "'313.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'0.9951"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.3946"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("D8").Value = "'0.3936"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'52.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.55%  "
$ws.Range("D10").Value = "'1.409"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").Value = "'0.9954"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'0.08606"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "'24.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'7.311"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "'0.00001340"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").Value = "'7.775"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "1.657.87"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "'95.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "'0.06955"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'20.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'6.995"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'0.9942"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'13.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "24.391.62"
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").Value = "'2.479"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.81%  "
$ws.Range("D26").Value = "'2.952"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.54%  "
$ws.Range("D27").Value = "'22.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'158.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").Value = "'142.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'5.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "'8.171"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "'2.563"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "1.829.71"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").Value = "'1.066"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.47%  "
$ws.Range("D35").Value = "'0.08250"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D40").Value = "'0.09276"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").Value = "'13.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("D42").Value = "'0.7767"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "'1.445"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'16.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").Value = "'0.7118"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("D46").Value = "'2.533"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "'4.135"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("D48").Value = "'0.9951"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "'0.08448"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'1.459"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.45%  "
$ws.Range("D51").Value = "'136.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.94%  "
